# Auto-generated Excel COM-interop script to apply scheduled price/profit refresh
# across the Kujata_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 900
$ws.Range("J12").Value = 900
$ws.Range("L12").Value = 900
$ws.Range("N12").Value = -1240
$ws.Range("H15").Value = 1199.64
$ws.Range("I15").Value = 1199.64
$ws.Range("K15").Value = 3598.92
$ws.Range("M15").Value = -3429.92
$ws.Range("H33").Value = 198.23077
$ws.Range("I33").Value = 147.90909
$ws.Range("J33").Value = 475
$ws.Range("K33").Value = 147.90909
$ws.Range("L33").Value = 475
$ws.Range("M33").Value = 81.09091000000001
$ws.Range("N33").Value = -933
$ws.Range("H40").Value = 1600.875
$ws.Range("I40").Value = 2749.5
$ws.Range("J40").Value = 1218
$ws.Range("K40").Value = 2749.5
$ws.Range("L40").Value = 1218
$ws.Range("M40").Value = -2574.5
$ws.Range("N40").Value = -1568
$ws.Range("H98").Value = 8020.9443
$ws.Range("I98").Value = 9225.267
$ws.Range("K98").Value = 9225.267
$ws.Range("M98").Value = -7727.267
$ws.Range("H112").Value = 2317.8572
$ws.Range("J112").Value = 2317.8572
$ws.Range("L112").Value = 6953.571599999999
$ws.Range("N112").Value = -9169.571599999999
$ws.Range("H116").Value = 2955.125
$ws.Range("I116").Value = 2941.889
$ws.Range("J116").Value = 2972.1428
$ws.Range("K116").Value = 2941.889
$ws.Range("L116").Value = 2972.1428
$ws.Range("M116").Value = 500.1109999999999
$ws.Range("N116").Value = -9856.1428
$ws.Range("H122").Value = 8020.9443
$ws.Range("I122").Value = 9225.267
$ws.Range("K122").Value = 27675.801
$ws.Range("M122").Value = -25225.801
$ws.Range("H127").Value = 2292.077
$ws.Range("I127").Value = 943.5
$ws.Range("J127").Value = 2537.2727
$ws.Range("K127").Value = 2830.5
$ws.Range("L127").Value = 7611.8181
$ws.Range("M127").Value = 2129.5
$ws.Range("N127").Value = -17531.8181
$ws.Range("H129").Value = 926.1111
$ws.Range("J129").Value = 990.3333
$ws.Range("L129").Value = 2970.9999
$ws.Range("N129").Value = -12970.9999
$ws.Range("H137").Value = 1401.35
$ws.Range("I137").Value = 870.1
$ws.Range("J137").Value = 1932.6
$ws.Range("K137").Value = 2610.3
$ws.Range("L137").Value = 5797.799999999999
$ws.Range("M137").Value = -60.30000000000018
$ws.Range("N137").Value = -10897.8
$ws.Range("H138").Value = 2859.1042
$ws.Range("I138").Value = 4898.6
$ws.Range("J138").Value = 2747.044
$ws.Range("K138").Value = 14695.8
$ws.Range("L138").Value = 8241.132
$ws.Range("M138").Value = -9555.800000000001
$ws.Range("N138").Value = -18521.132

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10654.14
$ws.Range("I32").Value = 7690.576
$ws.Range("J32").Value = 14918.78
$ws.Range("K32").Value = 7690.576
$ws.Range("L32").Value = 14918.78
$ws.Range("M32").Value = -7403.576
$ws.Range("N32").Value = -15492.78
$ws.Range("H88").Value = 1305.4
$ws.Range("I88").Value = 1430
$ws.Range("K88").Value = 1430
$ws.Range("M88").Value = -1024
$ws.Range("H91").Value = 1305.4
$ws.Range("I91").Value = 1430
$ws.Range("K91").Value = 1430
$ws.Range("M91").Value = -26

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4561.65
$ws.Range("I86").Value = 4729.4287
$ws.Range("K86").Value = 4729.4287
$ws.Range("M86").Value = -3606.4287
$ws.Range("H89").Value = 4561.65
$ws.Range("I89").Value = 4729.4287
$ws.Range("K89").Value = 23647.1435
$ws.Range("M89").Value = -18031.1435
$ws.Range("H99").Value = 55556764
$ws.Range("I99").Value = 71429630
$ws.Range("J99").Value = 1712.5
$ws.Range("K99").Value = 71429630
$ws.Range("L99").Value = 1712.5
$ws.Range("M99").Value = -71428132
$ws.Range("N99").Value = -4708.5
$ws.Range("H132").Value = 38786.668
$ws.Range("J132").Value = 38786.668
$ws.Range("L132").Value = 38786.668
$ws.Range("N132").Value = -48906.668

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 250.84616
$ws.Range("I7").Value = 217.88889
$ws.Range("J7").Value = 325
$ws.Range("K7").Value = 217.88889
$ws.Range("L7").Value = 325
$ws.Range("M7").Value = -104.88889
$ws.Range("N7").Value = -551
$ws.Range("H31").Value = 1790.3492
$ws.Range("I31").Value = 1736.2642
$ws.Range("J31").Value = 2077
$ws.Range("K31").Value = 1736.2642
$ws.Range("L31").Value = 2077
$ws.Range("M31").Value = -1441.2642
$ws.Range("N31").Value = -2667
$ws.Range("H34").Value = 1790.3492
$ws.Range("I34").Value = 1736.2642
$ws.Range("J34").Value = 2077
$ws.Range("K34").Value = 1736.2642
$ws.Range("L34").Value = 2077
$ws.Range("M34").Value = -1534.2642
$ws.Range("N34").Value = -2481
$ws.Range("H99").Value = 1755948.5
$ws.Range("I99").Value = 4387288
$ws.Range("K99").Value = 4387288
$ws.Range("M99").Value = -4385790
$ws.Range("H126").Value = 1755948.5
$ws.Range("I126").Value = 4387288
$ws.Range("K126").Value = 13161864
$ws.Range("M126").Value = -13159394
$ws.Range("H132").Value = 2159.9167
$ws.Range("I132").Value = 1759.0625
$ws.Range("J132").Value = 2961.625
$ws.Range("K132").Value = 5277.1875
$ws.Range("L132").Value = 8884.875
$ws.Range("M132").Value = -2747.1875
$ws.Range("N132").Value = -13944.875
$ws.Range("H141").Value = 339076.62
$ws.Range("J141").Value = 339076.62
$ws.Range("L141").Value = 339076.62
$ws.Range("N141").Value = -349436.62

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 54.5
$ws.Range("I2").Value = 54.5
$ws.Range("K2").Value = 327
$ws.Range("M2").Value = -214
$ws.Range("H68").Value = 1221.7142
$ws.Range("J68").Value = 1116.6666
$ws.Range("L68").Value = 3349.9998
$ws.Range("N68").Value = -4971.9998
$ws.Range("H71").Value = 1221.7142
$ws.Range("J71").Value = 1116.6666
$ws.Range("L71").Value = 10049.9994
$ws.Range("N71").Value = -18161.9994
$ws.Range("H132").Value = 1314
$ws.Range("I132").Value = 1245.4
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 11208.6
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -8678.6
$ws.Range("N132").Value = -23060
$ws.Range("H140").Value = 27066.805
$ws.Range("I140").Value = 55094.105
$ws.Range("J140").Value = 2861.4092
$ws.Range("K140").Value = 165282.315
$ws.Range("L140").Value = 8584.2276
$ws.Range("M140").Value = -160102.315
$ws.Range("N140").Value = -18944.2276

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3263.3333
$ws.Range("I80").Value = 1795
$ws.Range("K80").Value = 1795
$ws.Range("M80").Value = -797
$ws.Range("H83").Value = 3263.3333
$ws.Range("I83").Value = 1795
$ws.Range("K83").Value = 8975
$ws.Range("M83").Value = -3983
$ws.Range("H132").Value = 7848.6523
$ws.Range("I132").Value = 9053.352999999999
$ws.Range("K132").Value = 27160.059
$ws.Range("M132").Value = -24630.059

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1803.6428
$ws.Range("I82").Value = 1887.7273
$ws.Range("J82").Value = 1495.3334
$ws.Range("K82").Value = 1887.7273
$ws.Range("L82").Value = 1495.3334
$ws.Range("M82").Value = -1526.7273
$ws.Range("N82").Value = -2217.3334
$ws.Range("H85").Value = 1803.6428
$ws.Range("I85").Value = 1887.7273
$ws.Range("J85").Value = 1495.3334
$ws.Range("K85").Value = 1887.7273
$ws.Range("L85").Value = 1495.3334
$ws.Range("M85").Value = -639.7273
$ws.Range("N85").Value = -3991.3334
$ws.Range("H132").Value = 3382.7693
$ws.Range("I132").Value = 3398.2
$ws.Range("J132").Value = 3373.125
$ws.Range("K132").Value = 10194.6
$ws.Range("L132").Value = 10119.375
$ws.Range("M132").Value = -7664.599999999999
$ws.Range("N132").Value = -15179.375

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1868.2273
$ws.Range("I81").Value = 1600.5
$ws.Range("J81").Value = 1895
$ws.Range("K81").Value = 3201
$ws.Range("L81").Value = 3790
$ws.Range("M81").Value = -2140
$ws.Range("N81").Value = -5912
$ws.Range("H84").Value = 1868.2273
$ws.Range("I84").Value = 1600.5
$ws.Range("J84").Value = 1895
$ws.Range("K84").Value = 16005
$ws.Range("L84").Value = 18950
$ws.Range("M84").Value = -10701
$ws.Range("N84").Value = -29558
$ws.Range("H108").Value = 25999.334
$ws.Range("J108").Value = 25999.334
$ws.Range("L108").Value = 25999.334
$ws.Range("N108").Value = -33679.334
$ws.Range("H136").Value = 1650.3572
$ws.Range("I136").Value = 1283.3334
$ws.Range("J136").Value = 1925.625
$ws.Range("K136").Value = 3850.0002
$ws.Range("L136").Value = 5776.875
$ws.Range("M136").Value = -1300.0002
$ws.Range("N136").Value = -10876.875
